# Adds the "form" statistics columns (S, SC, QS, QSC, BTTS%, Blank%, CS%, FG%, FC%)
# to both the "home" and "away" sheets of the Bundesliga_Home_Away workbook.

$wb = $excel.ActiveWorkbook

$headers = @("S", "SC", "QS", "QSC", "BTTS%", "Blank%", "CS%", "FG%", "FC%")

# row, S, SC, QS, QSC, BTTS%, Blank%, CS%, FG%, FC%
$homeData = @(
    @(2,287,145,53,23,0.59,0.06,0.41,0.76,0.18),
    @(3,382,150,72,21,0.5,0.06,0.44,0.83,0.17),
    @(4,304,186,46,28,0.5,0.17,0.39,0.67,0.28),
    @(5,348,197,69,31,0.61,0.06,0.33,0.72,0.28),
    @(6,285,165,64,23,0.71,0.06,0.24,0.59,0.41),
    @(7,209,208,32,34,0.35,0.29,0.41,0.59,0.35),
    @(8,219,175,35,26,0.71,0.06,0.29,0.53,0.41),
    @(9,245,207,36,35,0.53,0.18,0.35,0.53,0.41),
    @(10,268,269,42,30,0.76,0,0.24,0.71,0.29),
    @(11,260,234,41,29,0.53,0.24,0.29,0.71,0.24),
    @(12,153,160,15,23,0.64,0.18,0.36,0.55,0.27),
    @(13,229,211,41,34,0.71,0.12,0.24,0.41,0.53),
    @(14,230,215,32,36,0.88,0.06,0.06,0.44,0.5600000000000001),
    @(15,195,248,27,39,0.82,0.12,0.06,0.47,0.53),
    @(16,206,211,31,34,0.71,0.12,0.18,0.47,0.53),
    @(17,228,165,27,29,0.62,0.25,0.12,0.38,0.62),
    @(18,258,198,21,42,0.47,0.47,0.18,0.35,0.53),
    @(19,169,156,19,31,0.58,0.42,0.08,0.17,0.75)
)

$awayData = @(
    @(2,256,187,38,27,0.47,0.18,0.47,0.53,0.35),
    @(3,268,171,52,21,0.6899999999999999,0.06,0.25,0.75,0.25),
    @(4,220,203,45,29,0.65,0.24,0.12,0.47,0.53),
    @(5,233,235,46,36,0.75,0.06,0.25,0.62,0.31),
    @(6,240,186,38,26,0.5600000000000001,0.19,0.25,0.44,0.5600000000000001),
    @(7,197,308,31,39,0.83,0.11,0.06,0.5,0.5),
    @(8,170,262,32,49,0.65,0.29,0.06,0.53,0.47),
    @(9,177,302,24,39,0.59,0.24,0.18,0.47,0.53),
    @(10,162,264,26,37,0.59,0.24,0.24,0.47,0.47),
    @(11,201,239,25,38,0.53,0.35,0.12,0.47,0.53),
    @(12,169,245,30,40,0.59,0.29,0.12,0.29,0.71),
    @(13,112,210,18,31,0.75,0.17,0.08,0.33,0.67),
    @(14,199,279,27,53,0.59,0.35,0.18,0.24,0.65),
    @(15,195,249,23,38,0.41,0.47,0.18,0.24,0.71),
    @(16,201,245,29,39,0.61,0.33,0.17,0.28,0.61),
    @(17,198,266,24,46,0.65,0.35,0,0.24,0.76),
    @(18,123,202,10,41,0.64,0.36,0.09,0.27,0.64),
    @(19,184,321,35,53,0.76,0.24,0.12,0.35,0.53)
)

function Fill-Sheet($ws, $data) {
    # Header row (J1:R1) - same bold/bordered/centered style as the other headers.
    for ($i = 0; $i -lt $headers.Length; $i++) {
        $ws.Cells.Item(1, 10 + $i).Value = $headers[$i]
    }
    $ws.Range("A1").Copy()
    $ws.Range("J1:R1").PasteSpecial(-4122)  # xlPasteFormats

    # Data rows (J2:R19)
    foreach ($row in $data) {
        $r = $row[0]
        for ($i = 1; $i -lt $row.Length; $i++) {
            $ws.Cells.Item($r, 9 + $i).Value = $row[$i]
        }
    }
}

$wsHome = $wb.Worksheets.Item("home")
Fill-Sheet $wsHome $homeData

$wsAway = $wb.Worksheets.Item("away")
Fill-Sheet $wsAway $awayData
